$d = $word.ActiveDocument

# Locate the last paragraph in the document (the "Bibliografia" body
# text paragraph) so the new content can be appended right before the
# final section properties, matching the target diff.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# --- New paragraph #1: "Requisitos" heading (Heading2 style) ---
$lastPara.Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$headingPara = $d.Paragraphs.Item($count)
$headingPara.Style = "Heading2"
$headingPara.Range.Text = "Requisitos"

# --- New paragraph #2: requirement bullet item (ListBullet style) ---
$headingPara.Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$reqPara = $d.Paragraphs.Item($count)
$reqPara.Style = "ListBullet"

# Write the text followed by a placeholder marker, then swap the
# marker for a manual line break ("^l") via Find/Replace so the break
# ends up inside the same run as the text, as in the target markup.
$reqPara.Range.Text = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)#MARK#"
$find = $reqPara.Range.Find
$find.Execute("#MARK#", $true, $false, $false, $false, $false, $true, 1, $false, "^l", 2)
